# Apply the commit: insert a new weekly price record for Tomate at row 181,
# pushing all subsequent rows (181-282) down by one (to 182-283).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 181; this shifts existing rows 181..282 to 182..283
# and copies formatting (e.g. the date style on column D) from the row above.
$ws.Rows.Item(181).Insert()

# Populate the newly inserted row 181 with the new record's data.
$ws.Cells.Item(181, 1).Value  = 7
$ws.Cells.Item(181, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(181, 3).Value  = "Ñuble"
$ws.Cells.Item(181, 4).Value  = 44460
$ws.Cells.Item(181, 5).Value  = 16
$ws.Cells.Item(181, 6).Value  = 100112020
$ws.Cells.Item(181, 7).Value  = "Tomate"
$ws.Cells.Item(181, 8).Value  = "Larga vida"
$ws.Cells.Item(181, 9).Value  = "Primera"
$ws.Cells.Item(181, 10).Value = 600
$ws.Cells.Item(181, 11).Value = 7000
$ws.Cells.Item(181, 12).Value = 7500
$ws.Cells.Item(181, 13).Value = 7250
$ws.Cells.Item(181, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(181, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(181, 16).Value = 725
$ws.Cells.Item(181, 17).Value = 10
$ws.Cells.Item(181, 18).Value = "Hortaliza"
